# Apply the crypto price/volume updates from the "Updated cryptos list" commit.
# Each target cell value is written with a literal leading apostrophe so Excel
# keeps it as text (matching the original inlineStr cells) instead of silently
# re-interpreting numeric-looking strings (e.g. "0.999", "52.236.39",
# "  +0.85%  ") as numbers and losing exact formatting such as trailing zeros,
# multiple "." separators, or padding spaces.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''52.236.39'
$ws.Range('E2').Value = '''  +0.85%  '
$ws.Range('D3').Value = '''2.882.76'
$ws.Range('E3').Value = '''  +3.43%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '''  -0.12%  '
$ws.Range('D5').Value = '''352.98'
$ws.Range('E5').Value = '''  +0.20%  '
$ws.Range('D6').Value = '''113.35'
$ws.Range('E6').Value = '''  +3.91%  '
$ws.Range('D7').Value = '''0.560'
$ws.Range('E7').Value = '''  +1.41%  '
$ws.Range('E8').Value = '''  +0.04%  '
$ws.Range('E9').Value = '''  +3.60%  '
$ws.Range('D10').Value = '''40.57'
$ws.Range('E10').Value = '''  +1.65%  '
$ws.Range('E11').Value = '''  -0.62%  '
$ws.Range('D12').Value = '''0.0855'
$ws.Range('E12').Value = '''  +1.98%  '
$ws.Range('D13').Value = '''20.19'
$ws.Range('E13').Value = '''  +0.06%  '
$ws.Range('E14').Value = '''  +2.36%  '
$ws.Range('E15').Value = '''  +3.07%  '
$ws.Range('E16').Value = '''  +4.34%  '
$ws.Range('D17').Value = '''0.995'
$ws.Range('E17').Value = '''  +7.35%  '
$ws.Range('D18').Value = '''52.210.34'
$ws.Range('E18').Value = '''  +0.87%  '
$ws.Range('E19').Value = '''  +8.11%  '
$ws.Range('E20').Value = '''  -1.16%  '
$ws.Range('E21').Value = '''  +3.31%  '
$ws.Range('D22').Value = '''0.0₃0979'
$ws.Range('E22').Value = '''  +1.37%  '
$ws.Range('E23').Value = '''  +1.21%  '
$ws.Range('D24').Value = '''270.78'
$ws.Range('E24').Value = '''  +1.39%  '
$ws.Range('E25').Value = '''  +1.79%  '
$ws.Range('D26').Value = '''26.64'
$ws.Range('E26').Value = '''  +1.80%  '
$ws.Range('E27').Value = '''  +0.02%  '
$ws.Range('E28').Value = '''  +1.73%  '
$ws.Range('B29').Value = '''Cosmos'
$ws.Range('C29').Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.61'
$ws.Range('E29').Value = '''  +3.75%  '
$ws.Range('B30').Value = '''InjectiveProtocol'
$ws.Range('C30').Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''39.06'
$ws.Range('E30').Value = '''  +5.83%  '
$ws.Range('E31').Value = '''  +2.33%  '
$ws.Range('E32').Value = '''  +1.72%  '
$ws.Range('D33').Value = '''0.0455'
$ws.Range('E33').Value = '''  +0.47%  '
$ws.Range('D34').Value = '''0.0899'
$ws.Range('E34').Value = '''  +8.17%  '
$ws.Range('E35').Value = '''  +1.83%  '
$ws.Range('B36').Value = '''Toncoin'
$ws.Range('C36').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').Value = '''1.94'
$ws.Range('E36').Value = '''  -12.96%  '
$ws.Range('B37').Value = '''FirstDigitalUSD'
$ws.Range('C37').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '''  -0.05%  '
$ws.Range('D38').Value = '''19.08'
$ws.Range('E38').Value = '''  +3.02%  '
$ws.Range('D39').Value = '''3.32'
$ws.Range('E39').Value = '''  +5.56%  '
$ws.Range('D40').Value = '''2.05'
$ws.Range('E40').Value = '''  +3.71%  '
$ws.Range('E41').Value = '''  +1.86%  '
$ws.Range('E42').Value = '''  +1.73%  '
$ws.Range('B43').Value = '''EnergySwap'
$ws.Range('C43').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''22.55'
$ws.Range('E43').Value = '''  +2.46%  '
$ws.Range('B44').Value = '''Monero'
$ws.Range('C44').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''122.56'
$ws.Range('E44').Value = '''  +1.50%  '
$ws.Range('E45').Value = '''  +1.83%  '
$ws.Range('E46').Value = '''  +7.76%  '
$ws.Range('D47').Value = '''2.182.07'
$ws.Range('E47').Value = '''  +2.76%  '
$ws.Range('E48').Value = '''  +7.14%  '
$ws.Range('E49').Value = '''  +17.57%  '
$ws.Range('D50').Value = '''0.959'
$ws.Range('E50').Value = '''  +5.46%  '
$ws.Range('D51').Value = '''0.0323'
$ws.Range('E51').Value = '''  +13.81%  '
